$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

# Rows 3-46: only the AttributeValues text (column A) changes; Size/Utility/UtilityDiff stay the same
$ws.Range("A3").Value = '{''Student'': np.int64(1), ''UndergradMajor'': np.int64(2)}'
$ws.Range("A4").Value = '{''Student'': np.int64(1), ''DevType'': np.int64(2)}'
$ws.Range("A5").Value = '{''Student'': np.int64(1), ''Gender'': np.int64(1)}'
$ws.Range("A6").Value = '{''Student'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Range("A9").Value = '{''UndergradMajor'': np.int64(2), ''Gender'': np.int64(1)}'
$ws.Range("A11").Value = '{''HDI'': np.int64(1), ''UndergradMajor'': np.int64(2)}'
$ws.Range("A12").Value = '{''DevType'': np.int64(2), ''Gender'': np.int64(1)}'
$ws.Range("A27").Value = '{''Student'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A31").Value = '{''Dependents'': np.int64(2), ''Hobby'': np.int64(1)}'
$ws.Range("A33").Value = '{''Dependents'': np.int64(2), ''Gender'': np.int64(1)}'
$ws.Range("A34").Value = '{''HDI'': np.int64(1), ''Gender'': np.int64(1)}'
$ws.Range("A35").Value = '{''Dependents'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Range("A36").Value = '{''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Range("A37").Value = '{''Dependents'': np.int64(2), ''HDI'': np.int64(1)}'
$ws.Range("A38").Value = '{''Gender'': np.int64(1), ''Student'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A39").Value = '{''Student'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A40").Value = '{''HDI'': np.int64(1), ''Student'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A41").Value = '{''Gender'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A42").Value = '{''Dependents'': np.int64(2), ''Gender'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A43").Value = '{''HDI'': np.int64(1), ''Gender'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A44").Value = '{''Dependents'': np.int64(2), ''SexualOrientation'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A45").Value = '{''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("A46").Value = '{''Student'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Gender'': np.int64(1)}'

# Rows 47-62: subgroup ordering rotated (old row 62 becomes new row 47; old rows 47-61 shift to 48-62)
$ws.Range("A47").Value = '{''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Gender'': np.int64(1), ''Hobby'': np.int64(1), ''Student'': np.int64(1)}'
$ws.Range("B47").Value = 15779
$ws.Range("C47").Value = 5234.178807957504
$ws.Range("D47").Value = -7729.070692964642
$ws.Range("A48").Value = '{''Dependents'': np.int64(2), ''Student'': np.int64(1), ''Gender'': np.int64(1)}'
$ws.Range("B48").Value = 15339
$ws.Range("C48").Value = 12445.08873558062
$ws.Range("D48").Value = -518.1607653415285
$ws.Range("A49").Value = '{''HDI'': np.int64(1), ''Student'': np.int64(1), ''Gender'': np.int64(1)}'
$ws.Range("B49").Value = 20079
$ws.Range("C49").Value = 8320.806949419126
$ws.Range("D49").Value = -4642.44255150302
$ws.Range("A50").Value = '{''Student'': np.int64(1), ''Dependents'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Range("B50").Value = 15301
$ws.Range("C50").Value = 12493.22924053899
$ws.Range("D50").Value = -470.0202603831513
$ws.Range("A51").Value = '{''Student'': np.int64(1), ''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Range("B51").Value = 20128
$ws.Range("C51").Value = 8058.492947675052
$ws.Range("D51").Value = -4904.756553247094
$ws.Range("A52").Value = '{''SexualOrientation'': np.int64(1), ''UndergradMajor'': np.int64(2), ''Gender'': np.int64(1)}'
$ws.Range("B52").Value = 16152
$ws.Range("C52").Value = 19730.73379293201
$ws.Range("D52").Value = 6767.48429200986
$ws.Range("A53").Value = '{''SexualOrientation'': np.int64(1), ''DevType'': np.int64(2), ''Gender'': np.int64(1)}'
$ws.Range("B53").Value = 15929
$ws.Range("C53").Value = 10271.32109283204
$ws.Range("D53").Value = -2691.928408090102
$ws.Range("A54").Value = '{''Dependents'': np.int64(2), ''SexualOrientation'': np.int64(1), ''Gender'': np.int64(1)}'
$ws.Range("B54").Value = 17563
$ws.Range("C54").Value = 12896.98458729801
$ws.Range("D54").Value = -66.26491362413253
$ws.Range("A55").Value = '{''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Gender'': np.int64(1)}'
$ws.Range("B55").Value = 21900
$ws.Range("C55").Value = 8524.583644570461
$ws.Range("D55").Value = -4438.665856351685
$ws.Range("A56").Value = '{''Dependents'': np.int64(2), ''HDI'': np.int64(1), ''Gender'': np.int64(1)}'
$ws.Range("B56").Value = 15786
$ws.Range("C56").Value = 12632.85656038428
$ws.Range("D56").Value = -330.392940537864
$ws.Range("A57").Value = '{''Dependents'': np.int64(2), ''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Range("B57").Value = 15673
$ws.Range("C57").Value = 12988.61598577243
$ws.Range("D57").Value = 25.36648485028491
$ws.Range("A58").Value = '{''Student'': np.int64(1), ''Gender'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("B58").Value = 18136
$ws.Range("C58").Value = 7283.726461989329
$ws.Range("D58").Value = -5679.523038932816
$ws.Range("A59").Value = '{''HDI'': np.int64(1), ''Gender'': np.int64(1), ''Student'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("B59").Value = 16648
$ws.Range("C59").Value = 6709.311194402443
$ws.Range("D59").Value = -6253.938306519703
$ws.Range("A60").Value = '{''Student'': np.int64(1), ''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("B60").Value = 16486
$ws.Range("C60").Value = 5663.715061677382
$ws.Range("D60").Value = -7299.534439244764
$ws.Range("A61").Value = '{''HDI'': np.int64(1), ''Gender'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Range("B61").Value = 18306
$ws.Range("C61").Value = 6562.835401713361
$ws.Range("D61").Value = -6400.414099208785
$ws.Range("A62").Value = '{''Student'': np.int64(1), ''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Gender'': np.int64(1)}'
$ws.Range("B62").Value = 19043
$ws.Range("C62").Value = 7319.305231565133
$ws.Range("D62").Value = -5643.944269357013
